# The sheet originally held yearly data for 2005年..2020年 in rows 2..17.
# The update drops the oldest five years (2005年-2009年) and appends a new
# row for 2021年, so the sheet now covers 2010年..2021年 in rows 2..13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the five oldest years (rows 2-6: 2005年-2009年). Everything below
# shifts up by 5 rows, so the row that held 2010年 (old row 7) becomes row 2,
# and the row that held 2020年 (old row 17) becomes row 12.
$ws.Rows("2:6").Delete()

# Append a new row 13 for 2021年. Copy row 12 (2020年) into row 13 first so
# the new row inherits the same cell formatting/style, then overwrite the
# values with the 2021年 figures.
$ws.Range("A12:J12").Copy($ws.Range("A13"))

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 660
$ws.Range("C13").Value = 118067
$ws.Range("D13").Value = 5803
$ws.Range("E13").Value = 3497
$ws.Range("F13").Value = 283
$ws.Range("G13").Value = 18388
$ws.Range("H13").Value = 21623
$ws.Range("I13").Value = 77579
$ws.Range("J13").Value = 15041
